$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("C2").Value = 11.38
$ws.Range("B3").Value = 8.619999999999999
$ws.Range("E4").Value = 10.17
$ws.Range("D5").Value = 9.83
$ws.Range("F5").Value = 10.23
$ws.Range("G5").Value = 9.4
$ws.Range("E6").Value = 9.77
$ws.Range("I6").Value = 9.43
$ws.Range("E7").Value = 10.6
$ws.Range("F9").Value = 10.57
$ws.Range("J9").Value = 13.75
$ws.Range("I10").Value = 6.25
